$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells store plain numeric-looking text as shared strings (type "s"),
# not numbers. Prefixing with a leading apostrophe forces Excel to keep the
# entry as text instead of silently converting it to a numeric value (which
# would also introduce floating point rounding noise).
$ws.Range("D11").Value = "'9.36"
$ws.Range("B33").Value = "'6.96"
$ws.Range("D33").Value = "'8.47"
$ws.Range("B36").Value = "'77.48"
$ws.Range("C36").Value = "'16.74"
$ws.Range("D36").Value = "'94.21"
